$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010 and 2010-18")

# ---------------------------------------------------------------------------
# Insert two new rows at 113/114. Excel will shift the existing blank rows
# 113/114 down to 115/116, and the data currently in rows 116/117/118 down
# to 118/119/120 - which is exactly the layout in the target workbook.
# ---------------------------------------------------------------------------
$ws.Rows("113:114").Insert()

# ---------------------------------------------------------------------------
# Row 113: new "CW3M C733" data row
# ---------------------------------------------------------------------------
$ws.Range("A113").Value = "CW3M C733"
$ws.Range("B113").Value = "Demo_Baseline_2010-18_C733"
$ws.Range("C113").Value = "2010-18"

$ws.Range("D113").NumberFormat = "0.00"
$ws.Range("D113").Value = 936.2841594444443

$ws.Range("E113").NumberFormat = "0.00"
$ws.Range("E113").Value = 1890.2624918888889

$ws.Range("F113").NumberFormat = "0.00"
$ws.Range("F113").Interior.ColorIndex = 6
$ws.Range("F113").Value = 0.59877088888888885

$ws.Range("G113").NumberFormat = "0.00"
$ws.Range("G113").Value = 270.41205844444437

$ws.Range("H113").NumberFormat = "0.00"
$ws.Range("H113").Interior.ColorIndex = 6
$ws.Range("H113").Value = 0

$ws.Range("I113").NumberFormat = "0.00"
$ws.Range("I113").Value = 7.3199817777777776

$ws.Range("J113").NumberFormat = "0.00"
$ws.Range("J113").Interior.ColorIndex = 6
$ws.Range("J113").Value = 0

$ws.Range("K113").NumberFormat = "0.00"
$ws.Range("K113").Value = 662.87599011111115

$ws.Range("L113").NumberFormat = "0.00"
$ws.Range("L113").Value = 80.365177222222229

$ws.Range("M113").NumberFormat = "0.00"
$ws.Range("M113").Value = 1422.4611409999998

$ws.Range("N113").NumberFormat = "0.00"
$ws.Range("N113").Value = 939.33221444444439

$ws.Range("O113").NumberFormat = "0"
$ws.Range("O113").Interior.ColorIndex = 6
$ws.Range("O113").Value = 4565.8835446666662

$ws.Range("P113").NumberFormat = "0"
$ws.Range("P113").Interior.ThemeColor = 6
$ws.Range("P113").Interior.TintAndShade = 0.79998168889431442
$ws.Range("P113").Interior.Color = 14083579
$ws.Range("P113").Value = 1017.8816121111109

$ws.Range("Q113").NumberFormat = "0.00"
$ws.Range("Q113").Value = 0.15706044444444447

$ws.Range("R113").NumberFormat = "0.000000"
$ws.Range("R113").Value = [double]"4.0444444444444593E-5"

$ws.Range("S113").Value = "has David Richey's new water rights data"

# ---------------------------------------------------------------------------
# Row 114: average-of-above formula row
# ---------------------------------------------------------------------------
$ws.Rows(114).RowHeight = 28.8

$ws.Range("A114").Value = "CW3M C733"
$ws.Range("B114").WrapText = $true
$ws.Range("B114").Value = "Demo_Baseline_2010-18_C733 + old water rights"
$ws.Range("C114").Value = "2010-18"

$ws.Range("D114").NumberFormat = "0.00"
$ws.Range("D114").Formula = "=AVERAGE(D105:D113)"

$ws.Range("E114:N114").NumberFormat = "0.00"
$ws.Range("O114:P114").NumberFormat = "0"
$ws.Range("Q114").NumberFormat = "0.00"
$ws.Range("R114").NumberFormat = "0.000000"
$ws.Range("E114:R114").Formula = "=AVERAGE(E105:E113)"

# ---------------------------------------------------------------------------
# Rows 115 & 116: blank spacer rows (formatting only, same pattern)
# ---------------------------------------------------------------------------
foreach ($r in 115, 116) {
    $ws.Range("B$r").WrapText = $true
    $ws.Range("D$r`:N$r").NumberFormat = "0.00"
    $ws.Range("O$r`:P$r").NumberFormat = "0"
    $ws.Range("O$r`:P$r").Interior.ColorIndex = 0
    $ws.Range("Q$r").NumberFormat = "0.00"
    $ws.Range("R$r").NumberFormat = "0.000000"
}

# ---------------------------------------------------------------------------
# View state: selection moves to O97 (single cell)
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("O97").Select()
